# Update After Round 7
# Update the match-day date serials in column A (rows 2-19) to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 45772
    3  = 45773
    4  = 45774
    5  = 45772
    6  = 45772
    7  = 45772
    8  = 45773
    9  = 45774
    10 = 45774
    11 = 45774
    12 = 45771
    13 = 45773
    14 = 45773
    15 = 45771
    16 = 45773
    17 = 45774
    18 = 45774
    19 = 45773
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $updates[$row]
}
